# Updated symbol list on Tue Dec 27 14:45:40 UTC 2022 with GitHub Actions
#
# Refresh the "Price" (D) and, for a few rows, the "Volume(1h)" (E) columns
# of the crypto price table with the latest scraped values. Price cells are
# stored as text (they can carry significant trailing zeros, e.g. "3.400"),
# so a leading apostrophe is used to keep Excel from re-typing the literal
# as a number when it is assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.39"
$ws.Range("D3").Value = "'24.01"
$ws.Range("D4").Value = "'5.402"
$ws.Range("D5").Value = "'0.05899"
$ws.Range("D6").Value = "'3.400"
$ws.Range("D7").Value = "'6.509"
$ws.Range("D8").Value = "'0.8114"
$ws.Range("D9").Value = "'0.9281"
$ws.Range("D10").Value = "'0.1418"
$ws.Range("D11").Value = "'0.07391"
$ws.Range("D12").Value = "'0.03096"
$ws.Range("D13").Value = "'0.03077"
$ws.Range("D14").Value = "'0.09333"
$ws.Range("D15").Value = "'3.870"
$ws.Range("D16").Value = "'0.001581"
$ws.Range("D17").Value = "'0.04739"

$ws.Range("D18").Value = "'0.0005969"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").Value = "'0.005923"
$ws.Range("D20").Value = "'0.001250"
$ws.Range("D21").Value = "'0.004732"
$ws.Range("D22").Value = "'0.00008815"
$ws.Range("D24").Value = "'2.158"

$ws.Range("D27").Value = "'0.0002657"
$ws.Range("E27").Value = "26UpBotsUBXT"

$ws.Range("D40").Value = "'0.03875"
$ws.Range("D41").Value = "'0.006401"
$ws.Range("D42").Value = "'0.1068"

$ws.Range("D43").Value = "'0.003091"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"

$ws.Range("D44").Value = "'0.008524"
$ws.Range("D45").Value = "'0.00005216"
$ws.Range("D47").Value = "'0.6721"

$ws.Range("D48").Value = "'0.001946"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
